$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: Mean ---
$ws.Range("A33").Value = "Mean"
$ws.Range("E33").Formula = "=AVERAGE(E6:E30)"
$ws.Range("F33").Formula = "=AVERAGE(F6:F30)"
$ws.Range("G33").Formula = "=AVERAGE(G6:G30)"
$ws.Range("I33").Formula = "=AVERAGE(I6:I30)"
$ws.Range("J33").Formula = "=AVERAGE(J6:J30)"
$ws.Range("K33").Formula = "=AVERAGE(K6:K30)"
$ws.Range("L33").Formula = "=AVERAGE(L6:L30)"
$ws.Range("M33").Formula = "=AVERAGE(M6:M30)"
$ws.Range("N33").Formula = "=AVERAGE(N6:N30)"

# --- Row 34: Max ---
$ws.Range("A34").Value = "Max"
$ws.Range("E34").Formula = "=MAX(E6:E30)"
$ws.Range("F34").Formula = "=MAX(F6:F30)"
$ws.Range("G34").Formula = "=MAX(G6:G30)"
$ws.Range("I34").Formula = "=MAX(I6:I30)"
$ws.Range("J34").Formula = "=MAX(J6:J30)"
$ws.Range("K34").Formula = "=MAX(K6:K30)"
$ws.Range("L34").Formula = "=MAX(L6:L30)"
$ws.Range("M34").Formula = "=MAX(M6:M30)"
$ws.Range("N34").Formula = "=MAX(N6:N30)"

# --- Row 35: Min ---
$ws.Range("A35").Value = "Min"
$ws.Range("E35").Formula = "=MIN(E6:E30)"
$ws.Range("F35").Formula = "=MIN(F6:F30)"
$ws.Range("G35").Formula = "=MIN(G6:G30)"
$ws.Range("I35").Formula = "=MIN(I6:I30)"
$ws.Range("J35").Formula = "=MIN(J6:J30)"
$ws.Range("K35").Formula = "=MIN(K6:K30)"
$ws.Range("L35").Formula = "=MIN(L6:L30)"
$ws.Range("M35").Formula = "=MIN(M6:M30)"
$ws.Range("N35").Formula = "=MIN(N6:N30)"

# --- Row heights to match the rest of the table ---
$ws.Range("A33:A35").RowHeight = 17.25

# --- Formatting: reuse the existing "centered + bordered" style (as used by
#     columns C/D/E/G in the data rows) for every value cell in the new rows,
#     and the plain bordered style (as used by column B) for column B. ---
$ws.Range("C6:D6").Copy() | Out-Null
$ws.Range("C33:D35").PasteSpecial(-4122) | Out-Null

$ws.Range("E6").Copy() | Out-Null
$ws.Range("E33:E35").PasteSpecial(-4122) | Out-Null

$ws.Range("G6").Copy() | Out-Null
$ws.Range("F33:G35").PasteSpecial(-4122) | Out-Null
$ws.Range("H33:H35").PasteSpecial(-4122) | Out-Null
$ws.Range("I33:N35").PasteSpecial(-4122) | Out-Null

$ws.Range("B6").Copy() | Out-Null
$ws.Range("B33:B35").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Copy() | Out-Null
$ws.Range("A33:A35").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- View state: leave the selection on H38, similar to the saved workbook ---
$ws.Range("A28").Select()
$ws.Range("H38").Select()
